# This workbook is a discrete-event simulation log (arrival/service times
# for a multi-server queue). The edit replaces the logged run with a new
# simulation run: the "Tiempo Ocio Server" (server idle time) column is
# dropped, the remaining columns (B:E) are recomputed with new values, and
# the trailing rows for clients 58-60 are removed (the new run only
# produced 57 client records instead of 60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Tiempo Ocio Server" column (F) -- header + all its data.
$ws.Range("F1").EntireColumn.Delete()

# The new simulation run only has 57 client rows (rows 2-58); remove the
# old trailing rows 59-61 (clients 58-60) so the sheet ends at row 58.
$ws.Range("59:61").EntireRow.Delete()

# Write the new simulation results (id Cliente, Tiempo Llegada,
# Tiempo Atencion, Tiempo Espera, id Server) for clients 1-57.
$data = New-Object 'object[,]' 57,5
$data[0,0] = 1; $data[0,1] = 18.82; $data[0,2] = 19; $data[0,3] = 0.18; $data[0,4] = 10
$data[1,0] = 2; $data[1,1] = 23.94; $data[1,2] = 24; $data[1,3] = 0.06; $data[1,4] = 5
$data[2,0] = 3; $data[2,1] = 28.4; $data[2,2] = 29; $data[2,3] = 0.6; $data[2,4] = 10
$data[3,0] = 4; $data[3,1] = 46.86; $data[3,2] = 47; $data[3,3] = 0.14; $data[3,4] = 8
$data[4,0] = 5; $data[4,1] = 48.23; $data[4,2] = 49; $data[4,3] = 0.77; $data[4,4] = 10
$data[5,0] = 6; $data[5,1] = 63.02; $data[5,2] = 64; $data[5,3] = 0.98; $data[5,4] = 5
$data[6,0] = 7; $data[6,1] = 83.66; $data[6,2] = 84; $data[6,3] = 0.34; $data[6,4] = 5
$data[7,0] = 8; $data[7,1] = 112.29; $data[7,2] = 113; $data[7,3] = 0.71; $data[7,4] = 4
$data[8,0] = 9; $data[8,1] = 97.29000000000001; $data[8,2] = 114; $data[8,3] = 16.71; $data[8,4] = 5
$data[9,0] = 10; $data[9,1] = 111.87; $data[9,2] = 115; $data[9,3] = 3.13; $data[9,4] = 6
$data[10,0] = 11; $data[10,1] = 109.91; $data[10,2] = 116; $data[10,3] = 6.09; $data[10,4] = 7
$data[11,0] = 12; $data[11,1] = 132.98; $data[11,2] = 133; $data[11,3] = 0.02; $data[11,4] = 4
$data[12,0] = 13; $data[12,1] = 126.33; $data[12,2] = 134; $data[12,3] = 7.67; $data[12,4] = 5
$data[13,0] = 14; $data[13,1] = 151.51; $data[13,2] = 152; $data[13,3] = 0.49; $data[13,4] = 3
$data[14,0] = 15; $data[14,1] = 212.86; $data[14,2] = 213; $data[14,3] = 0.14; $data[14,4] = 4
$data[15,0] = 16; $data[15,1] = 173.49; $data[15,2] = 214; $data[15,3] = 40.51; $data[15,4] = 5
$data[16,0] = 17; $data[16,1] = 168.61; $data[16,2] = 215; $data[16,3] = 46.39; $data[16,4] = 6
$data[17,0] = 18; $data[17,1] = 176.75; $data[17,2] = 216; $data[17,3] = 39.25; $data[17,4] = 7
$data[18,0] = 19; $data[18,1] = 233.62; $data[18,2] = 234; $data[18,3] = 0.38; $data[18,4] = 5
$data[19,0] = 20; $data[19,1] = 222.08; $data[19,2] = 235; $data[19,3] = 12.92; $data[19,4] = 6
$data[20,0] = 21; $data[20,1] = 213.17; $data[20,2] = 236; $data[20,3] = 22.83; $data[20,4] = 7
$data[21,0] = 22; $data[21,1] = 241.47; $data[21,2] = 242; $data[21,3] = 0.53; $data[21,4] = 3
$data[22,0] = 23; $data[22,1] = 234.91; $data[22,2] = 243; $data[22,3] = 8.09; $data[22,4] = 4
$data[23,0] = 24; $data[23,1] = 262.28; $data[23,2] = 263; $data[23,3] = 0.72; $data[23,4] = 4
$data[24,0] = 25; $data[24,1] = 270.9; $data[24,2] = 271; $data[24,3] = 0.1; $data[24,4] = 2
$data[25,0] = 26; $data[25,1] = 281.21; $data[25,2] = 282; $data[25,3] = 0.79; $data[25,4] = 3
$data[26,0] = 27; $data[26,1] = 331.57; $data[26,2] = 332; $data[26,3] = 0.43; $data[26,4] = 3
$data[27,0] = 28; $data[27,1] = 369.06; $data[27,2] = 370; $data[27,3] = 0.9399999999999999; $data[27,4] = 1
$data[28,0] = 29; $data[28,1] = 337.88; $data[28,2] = 371; $data[28,3] = 33.12; $data[28,4] = 2
$data[29,0] = 30; $data[29,1] = 298.28; $data[29,2] = 372; $data[29,3] = 73.72; $data[29,4] = 3
$data[30,0] = 31; $data[30,1] = 311.82; $data[30,2] = 373; $data[30,3] = 61.18; $data[30,4] = 4
$data[31,0] = 32; $data[31,1] = 337.78; $data[31,2] = 374; $data[31,3] = 36.22; $data[31,4] = 5
$data[32,0] = 33; $data[32,1] = 324.22; $data[32,2] = 375; $data[32,3] = 50.78; $data[32,4] = 6
$data[33,0] = 34; $data[33,1] = 346.35; $data[33,2] = 376; $data[33,3] = 29.65; $data[33,4] = 7
$data[34,0] = 35; $data[34,1] = 354.86; $data[34,2] = 377; $data[34,3] = 22.14; $data[34,4] = 8
$data[35,0] = 36; $data[35,1] = 450.69; $data[35,2] = 451; $data[35,3] = 0.31; $data[35,4] = 2
$data[36,0] = 37; $data[36,1] = 391.04; $data[36,2] = 452; $data[36,3] = 60.96; $data[36,4] = 3
$data[37,0] = 38; $data[37,1] = 424.7; $data[37,2] = 453; $data[37,3] = 28.3; $data[37,4] = 4
$data[38,0] = 39; $data[38,1] = 390.95; $data[38,2] = 454; $data[38,3] = 63.05; $data[38,4] = 5
$data[39,0] = 40; $data[39,1] = 405.69; $data[39,2] = 455; $data[39,3] = 49.31; $data[39,4] = 6
$data[40,0] = 41; $data[40,1] = 405.14; $data[40,2] = 456; $data[40,3] = 50.86; $data[40,4] = 7
$data[41,0] = 42; $data[41,1] = 530.91; $data[41,2] = 531; $data[41,3] = 0.09; $data[41,4] = 2
$data[42,0] = 43; $data[42,1] = 513.27; $data[42,2] = 532; $data[42,3] = 18.73; $data[42,4] = 3
$data[43,0] = 44; $data[43,1] = 451.75; $data[43,2] = 533; $data[43,3] = 81.25; $data[43,4] = 4
$data[44,0] = 45; $data[44,1] = 460.29; $data[44,2] = 534; $data[44,3] = 73.70999999999999; $data[44,4] = 5
$data[45,0] = 46; $data[45,1] = 504.56; $data[45,2] = 535; $data[45,3] = 30.44; $data[45,4] = 6
$data[46,0] = 47; $data[46,1] = 467.07; $data[46,2] = 536; $data[46,3] = 68.93000000000001; $data[46,4] = 7
$data[47,0] = 48; $data[47,1] = 506.97; $data[47,2] = 537; $data[47,3] = 30.03; $data[47,4] = 8
$data[48,0] = 49; $data[48,1] = 545; $data[48,2] = 545; $data[48,3] = 0; $data[48,4] = 6
$data[49,0] = 50; $data[49,1] = 553.8200000000001; $data[49,2] = 554; $data[49,3] = 0.18; $data[49,4] = 5
$data[50,0] = 51; $data[50,1] = 526.37; $data[50,2] = 555; $data[50,3] = 28.63; $data[50,4] = 6
$data[51,0] = 52; $data[51,1] = 541.5700000000001; $data[51,2] = 556; $data[51,3] = 14.43; $data[51,4] = 7
$data[52,0] = 53; $data[52,1] = 647.3099999999999; $data[52,2] = 648; $data[52,3] = 0.6899999999999999; $data[52,4] = 9
$data[53,0] = 54; $data[53,1] = 562.14; $data[53,2] = 649; $data[53,3] = 86.86; $data[53,4] = 10
$data[54,0] = 55; $data[54,1] = 558; $data[54,2] = 650; $data[54,3] = 92; $data[54,4] = 1
$data[55,0] = 56; $data[55,1] = 562.86; $data[55,2] = 651; $data[55,3] = 88.14; $data[55,4] = 2
$data[56,0] = 57; $data[56,1] = 584.89; $data[56,2] = 652; $data[56,3] = 67.11; $data[56,4] = 3

$ws.Range("A2:E58").Value = $data
